# "fix: change plan template"
# The Time Slot sheet's lunch break is shortened (11:45-12:30 instead of
# 11:30-13:00), shifting slots 5-8 earlier by 30 minutes and slot 4's end
# time out to 11:45. Also tidy up two cells that were still pointing at a
# now-redundant duplicate cell style (J1 on Plan, A9 on Time Slot) so they
# match the plain centered/bordered style used everywhere else, and move
# the saved selection on the Time Slot sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Time Slot sheet: updated break/slot boundaries
# ---------------------------------------------------------------------
$timeSlot = $wb.Worksheets.Item("Time Slot")

$timeSlot.Range("C5").Value = 0.48958333333333331   # row5 End   11:45
$timeSlot.Range("B6").Value = 0.52083333333333337   # row6 Start 12:30
$timeSlot.Range("C6").Value = 0.5625                # row6 End   13:30
$timeSlot.Range("B7").Value = 0.5625                # row7 Start 13:30
$timeSlot.Range("C7").Value = 0.60416666666666663   # row7 End   14:30
$timeSlot.Range("B8").Value = 0.60416666666666663   # row8 Start 14:30
$timeSlot.Range("C8").Value = 0.64583333333333337   # row8 End   15:30
$timeSlot.Range("B9").Value = 0.64583333333333337   # row9 Start 15:30
$timeSlot.Range("C9").Value = 0.6875                # row9 End   16:30

# A9 used a redundant one-off duplicate style (identical center+border
# look to the common style already used by A2:A8) - reapply the common
# formatting so it collapses back onto that shared style.
$timeSlot.Range("A9").Borders.LineStyle = 1
$timeSlot.Range("A9").HorizontalAlignment = -4108   # xlCenter

# Saved selection moved from G9 to C10
$timeSlot.Activate()
$timeSlot.Range("C10").Select()

# ---------------------------------------------------------------------
# Plan sheet: J1 had the same redundant duplicate style as Time Slot!A9
# ---------------------------------------------------------------------
$plan = $wb.Worksheets.Item("Plan")
$plan.Range("J1").Borders.LineStyle = 1
$plan.Range("J1").HorizontalAlignment = -4108       # xlCenter
